$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "72.374.50"
$ws.Range("E2").Value = "  +0.79%  "

# Row 3
$ws.Range("D3").Value = "2.716.30"
$ws.Range("E3").Value = "  +3.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'601.03"
$ws.Range("E5").Value = "  -0.95%  "

# Row 6
$ws.Range("D6").Value = "'176.59"
$ws.Range("E6").Value = "  -1.62%  "

# Row 8
$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("D9").Value = "2.715.76"
$ws.Range("E9").Value = "  +3.16%  "

# Row 10
$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  +0.71%  "

# Row 12
$ws.Range("D12").Value = "'0.356"
$ws.Range("E12").Value = "  +2.04%  "

# Row 13
$ws.Range("D13").Value = "'5.03"
$ws.Range("E13").Value = "  -0.29%  "

# Row 14
$ws.Range("D14").Value = "3.210.67"
$ws.Range("E14").Value = "  +2.43%  "

# Row 15
$ws.Range("E15").Value = "  -0.24%  "

# Row 16
$ws.Range("D16").Value = "72.080.68"
$ws.Range("E16").Value = "  +0.58%  "

# Row 17
$ws.Range("D17").Value = "'26.45"
$ws.Range("E17").Value = "  -0.44%  "

# Row 18
$ws.Range("D18").Value = "2.711.48"
$ws.Range("E18").Value = "  +2.95%  "

# Row 19
$ws.Range("D19").Value = "'12.34"
$ws.Range("E19").Value = "  +7.31%  "

# Row 20
$ws.Range("D20").Value = "'8.18"
$ws.Range("E20").Value = "  +2.63%  "

# Row 21
$ws.Range("D21").Value = "'375.07"
$ws.Range("E21").Value = "  -2.11%  "

# Row 22
$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  +0.46%  "

# Row 23
$ws.Range("E23").Value = "  +2.44%  "

# Row 24
$ws.Range("D24").Value = "'72.51"
$ws.Range("E24").Value = "  -0.32%  "

# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "'4.41"
$ws.Range("E25").Value = "  -1.23%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").Value = "'9.89"
$ws.Range("E27").Value = "  -0.50%  "

# Row 28
$ws.Range("D28").Value = "2.854.72"
$ws.Range("E28").Value = "  +3.17%  "

# Row 29
$ws.Range("D29").Value = "'0.993"
$ws.Range("E29").Value = "  -0.59%  "

# Row 30
$ws.Range("D30").Value = "'0.0000101"
$ws.Range("E30").Value = "  +4.15%  "

# Row 31
$ws.Range("D31").Value = "'8.18"
$ws.Range("E31").Value = "  +1.69%  "

# Row 32
$ws.Range("D32").Value = "'510.89"
$ws.Range("E32").Value = "  -6.38%  "

# Row 33
$ws.Range("D33").Value = "'1.32"
$ws.Range("E33").Value = "  -0.80%  "

# Row 34
$ws.Range("E34").Value = "  +0.09%  "

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("D36").Value = "'164.25"
$ws.Range("E36").Value = "  -1.07%  "

# Row 37
$ws.Range("D37").Value = "'19.76"
$ws.Range("E37").Value = "  +2.70%  "

# Row 38
$ws.Range("D38").Value = "'19.12"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39
$ws.Range("D39").Value = "'1.40"
$ws.Range("E39").Value = "  +0.12%  "

# Row 40
$ws.Range("D40").Value = "'0.110"
$ws.Range("E40").Value = "  -4.13%  "

# Row 41
$ws.Range("D41").Value = "'1.82"
$ws.Range("E41").Value = "  -2.57%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.10"
$ws.Range("E42").Value = "  +1.23%  "

# Row 43
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("D44").Value = "'2.59"
$ws.Range("E44").Value = "  -2.02%  "

# Row 45
$ws.Range("D45").Value = "'0.336"
$ws.Range("E45").Value = "  +1.15%  "

# Row 46
$ws.Range("D46").Value = "'157.15"
$ws.Range("E46").Value = "  +4.11%  "

# Row 47
$ws.Range("D47").Value = "'39.53"
$ws.Range("E47").Value = "  +0.72%  "

# Row 48
$ws.Range("D48").Value = "'0.568"
$ws.Range("E48").Value = "  +6.01%  "

# Row 49
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "'1.79"
$ws.Range("E49").Value = "  +5.98%  "

# Row 50
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'3.77"
$ws.Range("E50").Value = "  +3.20%  "

# Row 51
$ws.Range("E51").Value = "  +1.12%  "
